# Add a new "2021" column (column R) to the maternal mortality rate table,
# mirroring the existing "2020" column (Q): same formatting, one new year
# header and one new data value per region row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the whole Q3:Q14 band onto R3:R14 first, so every
# new cell (including the thin-border separator cell R3) picks up the same
# look (borders, font, number format, alignment) as its neighbour to the
# left before we overwrite the values that need to change.
$ws.Range("Q3:Q14").Copy()
$ws.Range("R3:R14").PasteSpecial(-4122)   # xlPasteFormats

# Header (row 4): the new year.
$ws.Range("R4").Value = 2021

# Data rows (5-14): the 2021 value for each region.
$ws.Range("R5").Value = 33.299999999999997
$ws.Range("R6").Value = 38.299999999999997
$ws.Range("R7").Value = 31.7
$ws.Range("R8").Value = 98.7
$ws.Range("R9").Value = 157.19999999999999
$ws.Range("R10").Value = 24.9
$ws.Range("R11").Value = 38.4
$ws.Range("R12").Value = 15.1
$ws.Range("R13").Value = 14.6
$ws.Range("R14").Value = 21.7

# Move the active selection/cursor to where it ended up after the edit.
$ws.Range("S6").Select()
